$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NegativeLoginTest")

# Insert a new column before column D (old D/E shift right to E/F),
# matching the width used by column C.
$cw = $ws.Columns.Item(3).ColumnWidth()
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).ColumnWidth = $cw

# New column D header + values ("Type" of negative-login scenario)
$ws.Range("D1").Value = "Type"
$ws.Range("D2").Value = "credentials"
$ws.Range("D3").Value = "credentials"
$ws.Range("D4").Value = "credentials"
$ws.Range("D5").Value = "fieldRequired"
$ws.Range("D6").Value = "fieldRequired"

# Replace the plain "wrongusername" values with a real-looking email and
# wire them up as mailto hyperlinks (matching B2/B6 treatment).
$ws.Range("B3").Value = "wrong@username.com"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:wrong@username.com")
$ws.Range("B3").Style = "Hyperlink"

$ws.Range("B4").Value = "wrong@username.com"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:wrong@username.com")
$ws.Range("B4").Style = "Hyperlink"

# Updated expected-message copy (now column E after the insert)
$ws.Range("E2").Value = "Error: No match for Email and/or Password"
$ws.Range("E3").Value = "Error: No match for Email and/or Password"
$ws.Range("E4").Value = "Error: No match for Email and/or Password"

# Cursor position cosmetics, matching the saved selection in the edit.
$ws.Range("E13").Select()
